# Regional Availability Factor workbook update
# - "updated 4.0 files and mdl": the hydrogen combustion turbine / hydrogen
#   combined cycle capacity-credit RAF values on the RAF-capacity sheet are
#   raised from 0.3 to 1, the About sheet's "last updated" date stamp moves
#   forward, and the workbook is left with the RAF-capacity sheet active
#   (scrolled/zoomed to show the rows that were just edited).

$wb = $excel.ActiveWorkbook

$wsAbout      = $wb.Worksheets.Item("About")
$wsGeneration = $wb.Worksheets.Item("RAF-generation")
$wsDemand     = $wb.Worksheets.Item("RAF-demand-altering-techs")
$wsCapacity   = $wb.Worksheets.Item("RAF-capacity")

# --- About sheet: bump the last-updated date stamp (2024-03-15 -> 2024-03-28) ---
$wsAbout.Range("C1").Value = 45379

# --- RAF-capacity sheet: the actual data edit -------------------------------
# hydrogen combustion turbine / hydrogen combined cycle capacity credit RAF
$wsCapacity.Range("B24").Value = 1
$wsCapacity.Range("B25").Value = 1

# widen column A a bit (matches the newly-added <cols> entry for this sheet)
$wsCapacity.Columns.Item(1).ColumnWidth = 28.1

# minor column-width retouches on the other RAF sheets (sub-pixel nudges)
$wsGeneration.Columns.Item(1).ColumnWidth = 33.6
$wsDemand.Columns.Item(1).ColumnWidth = 33.6
$wsDemand.Columns.Item(2).ColumnWidth = 18.3

# --- View / selection state --------------------------------------------------
# Make RAF-capacity the active sheet/tab (was RAF-generation before).
$wsCapacity.Activate()
$wsCapacity.Range("B25").Select()
$excel.ActiveWindow.Zoom = 80

$wsGeneration.Range("B3").Select()
$wsDemand.Range("D37").Select()

$wsCapacity.Activate()
